$d = $word.ActiveDocument

# --- Create the three new character styles -------------------------------

$styleGaNStyle = $d.Styles.Add("GaNStyle", 2)
$styleGaNStyle.Font.Name = "Calibri"
$styleGaNStyle.Font.Size = 14

$styleGaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$styleGaNParagraph.Font.Name = "Calibri"
$styleGaNParagraph.Font.Size = 10

$styleGaNLinks = $d.Styles.Add("GaNLinks", 2)
$styleGaNLinks.Font.Name = "Calibri"
$styleGaNLinks.Font.Bold = $true
$styleGaNLinks.Font.Color = 8388608
$styleGaNLinks.Font.Size = 9.5
$styleGaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Kampagnendaten..." heading run -------------

$kampagnendatenText = "Kampagnendaten 2022 für das Sternbild Cygnus: 10. bis 19. August, 9. bis 18. September, 8. bis 17. Oktober"
$range = $d.Content
$range.Find.Execute($kampagnendatenText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($range.Find.Found) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
    $range.Find.Execute($kampagnendatenText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Mach mit..." paragraph run ---------------

$machMitText = "Mach mit an einer weltweiten Kampagne, die schwächsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Cygnus am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel."
$range2 = $d.Content
$range2.Find.Execute($machMitText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($range2.Find.Found) {
    $range2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Die Schaubilder..." credit run ---------------

$schaubilderText = "Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$range3 = $d.Content
$range3.Find.Execute($schaubilderText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($range3.Find.Found) {
    $range3.Style = "GaNLinks"
}
